# "Use <formatversion> as suffix for table headers" — rename the
# "_old"/"_new" column-header suffixes to the concrete format versions
# FV2410 / FV2504, turn the header row into a proper Excel Table
# ("Table1" over A1:U62) and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- rename header row (A1:U1): "_old" -> "_FV2410", "_new" -> "_FV2504" ---
$ws.Range("A1").Value = "Segmentname_FV2410"
$ws.Range("B1").Value = "Segmentgruppe_FV2410"
$ws.Range("C1").Value = "Segment_FV2410"
$ws.Range("D1").Value = "Datenelement_FV2410"
$ws.Range("E1").Value = "Segment ID_FV2410"
$ws.Range("F1").Value = "Code_FV2410"
$ws.Range("G1").Value = "Qualifier_FV2410"
$ws.Range("H1").Value = "Beschreibung_FV2410"
$ws.Range("I1").Value = "Bedingungsausdruck_FV2410"
$ws.Range("J1").Value = "Bedingung_FV2410"
$ws.Range("K1").Value = "diff"
$ws.Range("L1").Value = "Segmentname_FV2504"
$ws.Range("M1").Value = "Segmentgruppe_FV2504"
$ws.Range("N1").Value = "Segment_FV2504"
$ws.Range("O1").Value = "Datenelement_FV2504"
$ws.Range("P1").Value = "Segment ID_FV2504"
$ws.Range("Q1").Value = "Code_FV2504"
$ws.Range("R1").Value = "Qualifier_FV2504"
$ws.Range("S1").Value = "Beschreibung_FV2504"
$ws.Range("T1").Value = "Bedingungsausdruck_FV2504"
$ws.Range("U1").Value = "Bedingung_FV2504"

# --- turn A1:U62 into an Excel Table (adds xl/tables/table1.xml +
#     tableParts wiring); column names are picked up from the just
#     renamed header row, so they already read "..._FV2410"/"..._FV2504" ---
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U62"), $null, 1)
$lo.Name = "Table1"
$lo.TableStyle = ""

# --- freeze the header row (row 1) ---
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
